$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 2;  F = 1;  G = 1 },
    @{ Row = 3;  F = 3;  G = 1 },
    @{ Row = 4;  F = 5;  G = 1 },
    @{ Row = 5;  F = 3;  G = 5 },
    @{ Row = 6;  F = 5;  G = 5 },
    @{ Row = 7;  F = 6;  G = 6 },
    @{ Row = 8;  F = 8;  G = 6 },
    @{ Row = 9;  F = 10; G = 6 },
    @{ Row = 10; F = 8;  G = 10 },
    @{ Row = 11; F = 10; G = 10 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = "U"
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

$ws.Range("D2:G11").Select()
